$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 17 (pushes old rows 17-19 down to 18-20).
# Excel copies formatting from the row above (row 16) into the new row,
# so C17/E17 inherit styles 3 and 15 respectively, matching the target.
$ws.Rows("17:17").Insert()

# The old row 19 (now row 20) had "Người Đề nghị" in column D.
# We need to move it to column E and put "Kế Toán" into the now-empty D.
# Copy D20's cell format into E20 first (PasteSpecial formats only, -4122 = xlPasteFormats),
# then move the text value across, and finally overwrite D20.
$ws.Range("D20").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E20").Value = $ws.Range("D20").Value2
$ws.Range("D20").Value = "Kế Toán"

# Fill in the newly inserted row 17 with the "amount in words" labels.
$ws.Range("C17").Value = "Số tiền bằng chữ:"
$ws.Range("E17").Value = "`${amountVND}"

# Match the updated active selection shown in the diff.
$ws.Range("E17").Select() | Out-Null
